$d = $word.ActiveDocument

# --- Change 1: Professional Summary paragraph -------------------------
# "...affecting all Black and Asian-American voters, developed geospatial ML..."
# -> "...affecting 50M voters, developed geospatial ML..."
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial ML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial ML", 2) | Out-Null

# --- Change 2: "Partner - Siege Analytics" bullet --------------------
# Split the single run into three runs so the new "50M" is bold
# (matching the formatting already used elsewhere in that bullet).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*race coding errors affecting all Black and Asian-American voters*") {
        $rng = $p.Range.Duplicate
        $rng.Find.Execute("all Black and Asian-American", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0) | Out-Null
        $rng.Text = "50M"
        $rng.Bold = 1
        $rng.Font.Color = 5258796   # wdColor BGR encoding of RGB 2C3E50
        break
    }
}

# --- Change 3: Move "Analytics Supervisor - GSD&M" block ---------------
# It currently sits between "Partner - Siege Analytics" and
# "Data Products Manager"; move it so it follows the
# "Data Products Manager" block instead.
$startPara = -1
$endPara = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startPara -eq -1 -and $t -like "*Analytics Supervisor - GSD&M*") {
        $startPara = $i
    }
    if ($startPara -ne -1 -and $t -like "*Advanced Statistical and ML techniques for segmentation*") {
        $endPara = $i
        break
    }
}

$blockStart = $d.Paragraphs.Item($startPara)
$blockEnd = $d.Paragraphs.Item($endPara)
$moveRng = $d.Range($blockStart.Range.Start, $blockEnd.Range.End)
$moveRng.Cut() | Out-Null

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*reducing processing time by*57%*") {
        $target = $p
    }
}
$pasteRng = $d.Range($target.Range.End, $target.Range.End)
$pasteRng.Paste() | Out-Null

# The cut/paste drops the paragraph style of the moved heading, so
# restore it explicitly.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Analytics Supervisor - GSD&M*") {
        $p.Style = "Heading 3"
        break
    }
}

# --- Change 4: Key Projects "Impact" line ------------------------------
# "...affecting all Black and Asian-American voters, improved electoral..."
# -> "...affecting 50M voters nationwide, improved electoral..."
$d.Content.Find.Execute(
    "affecting all Black and Asian-American voters, improved electoral prediction accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters nationwide, improved electoral prediction accuracy", 2) | Out-Null
